{"js": "// The activity guide repeats an intro sentence several times\n// (\"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148\n// prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Pegasus 2022: 8.\u201317. \u0159\u00edjna, 7.\u201316.\n// listopadu\"). Replace every occurrence with the updated wording that\n// drops \"Souhv\u011bzd\u00ed Pegasus 2022:\" from that spot and appends a new\n// sentence (plus a trailing date repeat), matching the source commit.\nconst oldText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Pegasus 2022: 8.\\u201317. \u0159\u00edjna, 7.\\u201316. listopadu\";\nconst newText =\n  \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 8.\\u201317. \u0159\u00edjna, 7.\\u201316. listopadu. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Pegasus.8.\\u201317. \u0159\u00edjna, 7.\\u201316. listopadu\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The activity guide repeats an intro sentence several times\n# (\"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148\n# prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Pegasus 2022: 8.\u201317. \u0159\u00edjna, 7.\u201316.\n# listopadu\"). Replace every occurrence with the updated wording that\n# drops \"Souhv\u011bzd\u00ed Pegasus 2022:\" from that spot and appends a new\n# sentence (plus a trailing date repeat), matching the source commit.\n$d = $word.ActiveDocument\n\n$old = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Souhv\u011bzd\u00ed Pegasus 2022: 8.\u201317. \u0159\u00edjna, 7.\u201316. listopadu\"\n$new = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od 8.\u201317. \u0159\u00edjna, 7.\u201316. listopadu. P\u0159i pozorov\u00e1n\u00ed pou\u017eijte hv\u011bzdy oblohy, kter\u00e9 zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Pegasus.8.\u201317. \u0159\u00edjna, 7.\u201316. listopadu\"\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$find.Forward = $true\n$find.Wrap = $wdFindContinue\n$find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null\n"}
